$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the merged "name" value from the existing code/diameter/category
# columns (A/B/C) before those columns disappear. ---
$code = $ws.Range("A2").Value2
$diameter = $ws.Range("B2").Value2
$category = $ws.Range("C2").Value2
$name = "$code-$diameter-$category"

# --- Remove the "diameter" and "category" columns (B, C). "code" (A) stays
# in place and gets repurposed into the new "name" column; quantity/state/
# warehouse/serials (old D/E/F/G) shift left into B/C/D/E. ---
$ws.Range("B1:C1").EntireColumn.Delete() | Out-Null

# Repurpose column A: "code" header/value -> "name" header/value.
$ws.Range("A1").Value = "name"
$ws.Range("A2").Value = $name

# Capitalize the "state" value: "mới" -> "Mới".
$ws.Range("C2").Value = "Mới"
